# Apply crypto price/volume updates scraped on Thu Mar 16 03:31:03 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.406.93'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.89%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.38%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.09%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.13%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3639'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.98%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.96'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.28%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3261'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.47%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.124'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.03%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.79%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.15%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.942'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.07%  '

# Row 14
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -8.10%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.599'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.45%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.652.82'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.47%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -8.14%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06615'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.64%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.06%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.12'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.90%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.932'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.35%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.67'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -9.81%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.44'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.45%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.404.80'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.475'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.11%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.351'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -15.94%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.06'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.70%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -9.04%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.837.59'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.42%  '

# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.190'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.25%  '

# Row 31
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '124.22'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.39%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.085'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.52%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.665'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -18.21%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08439'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.08%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.658'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -10.21%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.32'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -11.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.198'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.71%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06040'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -9.95%  '

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.222'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.35%  '

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02218'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -8.32%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2070'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.79%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.194'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -12.26%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.14%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5910'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -8.65%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.773'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.42%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.63'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -9.53%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5631'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.91%  '

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.06%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.944'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -9.32%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.72%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.70'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.74%  '
